# Apply odds/value updates to Sheet1, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 1.96
$ws.Range("N2").Value = 3.6
$ws.Range("Q2").Value = 1.91
$ws.Range("R2").Value = 1.34
$ws.Range("U2").Value = 2.02
$ws.Range("W2").Value = 2.04
$ws.Range("AE2").Value = 65
$ws.Range("AG2").Value = 11

# Row 3
$ws.Range("G3").Value = 2.74
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 4.1
$ws.Range("M3").Value = 1.08
$ws.Range("V3").Value = 1.33
$ws.Range("W3").Value = 1.59

# Row 4
$ws.Range("AB4").Value = 980

# Row 5
$ws.Range("G5").Value = 1.88
$ws.Range("H5").Value = 4.4
$ws.Range("I5").Value = 5.2
$ws.Range("J5").Value = 3.8
$ws.Range("M5").Value = 1.04
$ws.Range("T5").Value = 1.64
$ws.Range("W5").Value = 2.14
$ws.Range("Z5").Value = 42
$ws.Range("AA5").Value = 120
$ws.Range("AE5").Value = 60

# Row 6
$ws.Range("N6").Value = 2.8
$ws.Range("P6").Value = 1.58

# Row 7
$ws.Range("H7").Value = 3.25
$ws.Range("L7").Value = 1.28
$ws.Range("V7").Value = 1.37

# Row 8
$ws.Range("N8").Value = 6.2
$ws.Range("O8").Value = 1.12
$ws.Range("Y8").Value = 1000
$ws.Range("Z8").Value = 1000
$ws.Range("AB8").Value = 1000
$ws.Range("AE8").Value = 1000
$ws.Range("AF8").Value = 1000
$ws.Range("AI8").Value = 1000
$ws.Range("AK8").Value = 1000
$ws.Range("AL8").Value = 1000

# Row 9
$ws.Range("K9").Value = 4.8
$ws.Range("L9").Value = 1.19
$ws.Range("AD9").Value = 1000
$ws.Range("AO9").Value = 1000

# Row 10
$ws.Range("K10").Value = 5

# Row 11
$ws.Range("H11").Value = 3.8
$ws.Range("L11").Value = 1.43
$ws.Range("M11").Value = 1.08
$ws.Range("S11").Value = 3.4
$ws.Range("AC11").Value = 9

# Row 12
$ws.Range("I12").Value = 2.88

# Row 13
$ws.Range("I13").Value = 16
$ws.Range("Q13").Value = 1.69

# Row 14
$ws.Range("L14").Value = 1.47

$wb.Save()
